# "Start of Section 3"
#
# Mark the "Section 2" task (row 5) on the "Project 1 - dotnet project"
# sheet as finished: bump its Actual Time Spent to the full 3.5h estimate
# and flip its Status to "Done". All the other cells touched by the
# original diff (Summary sheet roll-up formulas, shared-string indices
# shifting down once "In Progress" is no longer used anywhere, cached
# formula values, etc.) are downstream consequences of this single data
# edit and are recomputed automatically on recalc/save.
#
# The selection is also moved to E6 (the Status cell of "Section 3"),
# matching the author having clicked there to start working on it next.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Project 1 - dotnet project")

$ws.Range("C5").Value = 3.5
$ws.Range("E5").Value = "Done"

$ws.Range("E6").Select() | Out-Null
